$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("C1").Value = "message"
$ws.Range("C2").Value = "Hello"

$ws.Activate()
$ws.Range("C2").Select()
